$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "33÷8=" "23÷8="
Replace-Text "91÷7=" "50÷9="
Replace-Text "96÷5=" "74÷3="
Replace-Text "74÷4=" "20÷3="
Replace-Text "24÷6=" "29÷9="

Replace-Text "46÷7=" "40÷2="
Replace-Text "80÷7=" "33÷9="
Replace-Text "39÷3=" "34÷9="
Replace-Text "11÷5=" "35÷6="
Replace-Text "67÷2=" "32÷9="

Replace-Text "24÷4=" "16÷2="
Replace-Text "85÷4=" "37÷3="
Replace-Text "64÷3=" "46÷4="
Replace-Text "78÷8=" "49÷8="
Replace-Text "40÷4=" "50÷8="

Replace-Text "38÷7=" "46÷9="
Replace-Text "23÷7=" "23÷9="
Replace-Text "88÷4=" "13÷2="
Replace-Text "69÷4=" "44÷6="
Replace-Text "44÷4=" "51÷7="

# The last row has "82÷8=" appearing twice (columns 1 and 2), which must
# become two different values, so address those cells directly rather
# than via a global Find/Replace.
$t = $d.Tables.Item(1)
$t.Cell(17, 1).Range.Text = "98÷2="
$t.Cell(17, 2).Range.Text = "75÷2="

Replace-Text "33÷7=" "26÷4="
Replace-Text "50÷5=" "78÷4="
Replace-Text "87÷6=" "80÷4="
